$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 157
$ws.Range("I2").Value = 385
$ws.Range("J2").Value = 1710
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 445
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 320
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 19
$ws.Range("S2").Value = 196
$ws.Range("T2").Value = 297
$ws.Range("U2").Value = 23
$ws.Range("V2").Value = 2597
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 2596
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 42
$ws.Range("AA2").Value = 23
